$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a literature-survey break note to the existing last row (row 38, col I) ---
$existing = $ws.Range("I38").Value2
$ws.Range("I38").Value = $existing + "`nב-30.7 עשיתי הפסקה מתודית בשביל סקירות ספרות"
$ws.Rows.Item(38).RowHeight = 43.5

# --- New row 39: literature survey on "ensemble learning" ---
$ws.Range("A39").Value = "נספחים"
$ws.Range("B39").Value = "תאורטי"
$ws.Range("C39").Value = "סקר ספרות"
$ws.Range("I39").Value = "חפיפה גדולה עם הפרק על עצים"
$ws.Range("D39").Value = "לקרוא, לחקור, לסכם, להשוות, לשאול שאלות, למצוא תשובות על ensemble learning"
$ws.Range("D39").WrapText = $true
$ws.Range("E39").Value = 44042.367361111108
$ws.Range("F39").Value = 44042.462500000001
$ws.Range("G39").Value = 0.25
$ws.Range("H39").Value = 0.25
$ws.Rows.Item(39).RowHeight = 29

# --- New row 40: literature survey on "unbalanced dataset" ---
$ws.Range("A40").Value = "נספחים"
$ws.Range("B40").Value = "תאורטי"
$ws.Range("C40").Value = "סקר ספרות"
$ws.Range("D40").Value = "לקרוא, לחקור, לסכם, להשוות, לשאול שאלות, למצוא תשובות על unbalanced dataset"
$ws.Range("D40").WrapText = $true
$ws.Range("E40").Value = 44042.502083333333
$ws.Range("F40").Value = 44042.64166666667
$ws.Range("G40").Value = 0.5
$ws.Range("H40").Value = 0.25
$ws.Range("I40").Value = "לכוון לנושאים יותר ספציפיים כי תמיד מגיעים לאותן 2.5 שיטות לטיפול בבעיה.`nרלוונטי לכל הסקירות: להוסיף אותן במקומות ספציפייים בהכשרה (בשלב הזה זה מרגיש מאוחר מדי)"
$ws.Range("I40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 43.5

# --- Update the view/selection state to match the post-edit scroll position ---
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("A37") } catch {}
$ws.Range("G43").Select() | Out-Null
